$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as text so that values such as
# "240.63" are not auto-converted to numbers by Excel, matching the inline
# string (t="inlineStr") representation used in the workbook.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "41.950.70"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "2.211.15"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "240.63"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "73.11"
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "42.29"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "0.0951"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "7.06"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "2.544.82"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "14.14"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").Value = "2.189.41"
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("D18").Value = "41.869.56"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("E19").Value = "  +9.00%  "
$ws.Range("D20").Value = "72.46"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").Value = "10.21"
$ws.Range("E22").Value = "  +16.79%  "
$ws.Range("D23").Value = "228.98"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("E24").Value = "  -7.03%  "
$ws.Range("D25").Value = "11.65"
$ws.Range("E25").Value = "  +3.34%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").Value = "167.08"
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").Value = "20.54"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "5.64"
$ws.Range("E32").Value = "  +7.90%  "
$ws.Range("D33").Value = "0.0786"
$ws.Range("E33").Value = "  -3.84%  "
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").Value = "28.74"
$ws.Range("E35").Value = "  -5.11%  "
$ws.Range("D36").Value = "0.110"
$ws.Range("E36").Value = "  -7.85%  "
$ws.Range("E37").Value = "  -4.65%  "
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("D39").Value = "13.16"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("D41").Value = "64.75"
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("D44").Value = "8.69"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").Value = "103.64"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("E47").Value = "  +4.94%  "
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "2.418.06"
$ws.Range("E51").Value = "  -2.17%  "

# Restore the default cell style on column D (the text number format was only
# needed transiently to prevent automatic numeric conversion while assigning).
$priceRange.Style = "Normal"

